# ---------------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# Two independent changes were made to the deck:
#
#  1. On slide 5, the table (graphicFrame "Google Shape;122;p17") had its
#     table style switched from the default "Table_0" style
#     ({93899FF5-CEF7-4622-89B1-81BEEECD4853}) to
#     {096D1FDF-31CD-48FD-93F5-0D1E08324684}.
#
#  2. The presentation's design/theme was changed back from the custom
#     "Integral" (Red Violet) look to the plain default "Office Theme"
#     colors. (The slide master / presentation theme part keeps its file
#     identity; only the 12 theme colors change.)
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style on slide 5's table (shape #2: Shape1=placeholder title,
#    Shape2=graphicFrame/table, Shape3=textbox).
# ---------------------------------------------------------------------------

$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{096D1FDF-31CD-48FD-93F5-0D1E08324684}")

# ---------------------------------------------------------------------------
# 2) Restore the classic "Office Theme" color scheme on the slide master's
#    theme (the same theme object backs Application's active design).
# ---------------------------------------------------------------------------

function Convert-HexToBgr($hex) {
    # PowerPoint's RGB color values are stored as 0x00BBGGRR (the same
    # packing VBA's RGB() function produces), so swap a "RRGGBB" hex
    # string into that integer encoding.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$colorScheme = $p.SlideMaster.ColorScheme

# Order matches MsoThemeColorSchemeIndex / the <a:clrScheme> child order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$colorScheme.Colors(1).RGB  = Convert-HexToBgr "000000"  # dk1
$colorScheme.Colors(2).RGB  = Convert-HexToBgr "FFFFFF"  # lt1
$colorScheme.Colors(3).RGB  = Convert-HexToBgr "44546A"  # dk2
$colorScheme.Colors(4).RGB  = Convert-HexToBgr "E7E6E6"  # lt2
$colorScheme.Colors(5).RGB  = Convert-HexToBgr "5B9BD5"  # accent1
$colorScheme.Colors(6).RGB  = Convert-HexToBgr "ED7D31"  # accent2
$colorScheme.Colors(7).RGB  = Convert-HexToBgr "A5A5A5"  # accent3
$colorScheme.Colors(8).RGB  = Convert-HexToBgr "FFC000"  # accent4
$colorScheme.Colors(9).RGB  = Convert-HexToBgr "4472C4"  # accent5
$colorScheme.Colors(10).RGB = Convert-HexToBgr "70AD47"  # accent6
$colorScheme.Colors(11).RGB = Convert-HexToBgr "0563C1"  # hlink
$colorScheme.Colors(12).RGB = Convert-HexToBgr "954F72"  # folHlink
